# Update Name of Algo
# Apply corrected KNN-imputed values to specific cells in Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.934
$ws.Range("C8").Value = -12.672
$ws.Range("A12").Value = -21.882
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -12.049
$ws.Range("C22").Value = -12.929
